# Applies the benchSuite-dacapo_gc-shenandoahGC_app-lusearch_heap-4G.docx
# stats-table refresh described by the commit:
#   "Fixed README.md stats and docx preparation for all DaCapo -
#    JDK 17 - Shenandoah GC tests"
#
# The document is a single-column table. Most rows just get their
# numeric text value swapped for the refreshed measurement; the three
# "summary" rows near the bottom (which previously carried a whole
# tab-separated record) collapse down to a single value matching the
# corresponding row near the top of the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Top summary rows -----------------------------------------------
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "391"

# --- Per-iteration timing rows ---------------------------------------
$t.Cell(6, 1).Range.Text  = "0.00052"
$t.Cell(7, 1).Range.Text  = "0.00012"
$t.Cell(8, 1).Range.Text  = "0.00003"
$t.Cell(9, 1).Range.Text  = "0.00021"
$t.Cell(10, 1).Range.Text = "0.00023"
$t.Cell(11, 1).Range.Text = "0.00025"
$t.Cell(12, 1).Range.Text = "0.05393"

# --- Bottom rows: collapse the old tab-delimited records down to a
#     single refreshed value -------------------------------------------
$t.Cell(44, 1).Range.Text = "99.93"
$t.Cell(45, 1).Range.Text = "0.05"
$t.Cell(46, 1).Range.Text = "72"
